# Add new column R ("l1") to the repeat_p39 sheet.
# Mirrors the commit: a new repeating-group indicator column "l1" is added
# right after the existing "c1" column, with a 0/1 indicator value (0) for
# every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell R1: same bold/border/centered style as the other header cells
# (copy formatting from the neighboring header cell Q1), with the text "l1".
$ws.Range("Q1").Copy() | Out-Null
$ws.Range("R1").PasteSpecial(-4122) | Out-Null
$ws.Range("R1").Value = "l1"

# Data cells R2:R74: numeric indicator value 0, matching the P/Q columns.
for ($r = 2; $r -le 74; $r++) {
    $ws.Cells.Item($r, 18).Value = 0
}
